# TestData.xlsx - rename worksheets (per commit: "clread build and updatedreadme")
#   FileData -> SignUpTest
#   EditData -> SearchItem
$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("FileData").Name = "SignUpTest"
$wb.Worksheets.Item("EditData").Name = "SearchItem"
